$d = $word.ActiveDocument

$replacements = @(
    @{old="70×13=910"; new="71×47=3337"},
    @{old="66×41=2706"; new="95×22=2090"},
    @{old="60×90=5400"; new="52×95=4940"},
    @{old="79×16=1264"; new="74×22=1628"},
    @{old="13×88=1144"; new="29×65=1885"},
    @{old="71×62=4402"; new="47×12=564"},
    @{old="30×32=960"; new="83×96=7968"},
    @{old="30×56=1680"; new="23×84=1932"},
    @{old="94×18=1692"; new="28×83=2324"},
    @{old="32×53=1696"; new="28×55=1540"},
    @{old="58×82=4756"; new="78×94=7332"},
    @{old="68×55=3740"; new="65×19=1235"},
    @{old="19×83=1577"; new="96×82=7872"},
    @{old="68×31=2108"; new="14×77=1078"},
    @{old="32×44=1408"; new="44×88=3872"},
    @{old="30×91=2730"; new="38×90=3420"},
    @{old="17×89=1513"; new="93×71=6603"},
    @{old="25×91=2275"; new="71×48=3408"},
    @{old="61×67=4087"; new="82×13=1066"},
    @{old="88×70=6160"; new="64×60=3840"},
    @{old="93×40=3720"; new="45×17=765"},
    @{old="21×15=315"; new="83×28=2324"},
    @{old="73×74=5402"; new="17×22=374"},
    @{old="81×98=7938"; new="63×45=2835"},
    @{old="42×46=1932"; new="18×83=1494"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $found = $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $($r.old)"
    }
}

$d.Save()
